# Applies the diff: adds two reference-note paragraphs / sentences,
# splits one long paragraph into several, and inserts a Word
# "citation" content control (Raw17 / Rawal et al., 2017) plus a
# couple of lastRenderedPageBreak markers that shift around as the
# content grows.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------
# 1) Insert a new paragraph "Snee (2015) misconception..." right
#    before the "A recent outbreak" paragraph. The original
#    paragraph (with ind firstLine=720) is rewritten in place so
#    that it now holds the new "Snee" sentence, and the old
#    "A recent outbreak..." text is re-emitted as a fresh paragraph
#    immediately after it (matching the diff, which keeps the
#    original <w:p> for the Snee text and adds a brand-new <w:p> for
#    the outbreak paragraph).
# ---------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("A recent outbreak of food poisoning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$outbreakPara = $find.Parent.Paragraphs(1)
$r = $outbreakPara.Range
$r.Collapse(1)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Snee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (2015) misconception that lots of data + analysis = magic.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>A recent outbreak of food poisoning at some locations has damaged the brand&#8217;s image and caused a significant decrease in sales.  The leadership team wants to restore consumer confidence by operationalizing their data lake to answer targeted questions about the incident.  Which sites are likely to have an outbreak next?  Are food handling procedures being followed?  Who should promotional material target?</w:t></w:r></w:p>
'@
$r.InsertXML($xml)

# ---------------------------------------------------------------
# 2) The (now empty) paragraph that used to hold only a <w:tab/>
#    right after the outbreak paragraph gains a lastRenderedPageBreak
#    marker ahead of the tab (page break shifted earlier).
# ---------------------------------------------------------------
$find.Execute("Section II: Collecting and Enhancement", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sectionIIPara = $find.Parent.Paragraphs(1)
$tabPara = $sectionIIPara.Previous()
$tr = $tabPara.Range
$tr.Collapse(1)
$xmlTab = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:tab/></w:r></w:p>'
$tr.InsertXML($xmlTab)

# ---------------------------------------------------------------
# 3) Remove the lastRenderedPageBreak that used to sit in front of
#    "Section II: Collecting and Enhancement" (it moved to the
#    paragraph above in step 2).
# ---------------------------------------------------------------
$find.Execute("Section II: Collecting and Enhancement", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sectionIIPara = $find.Parent.Paragraphs(1)
$sr = $sectionIIPara.Range
$sr.Collapse(1)
$xmlSection2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Section II: Collecting and Enhancement</w:t></w:r></w:p>'
$sr.InsertXML($xmlSection2)

Write-Host "step1-3 done"

# ---------------------------------------------------------------
# 4) Split the "Before analysis can begin..." paragraph into four:
#      a) tab + "According to Gibert et al. (2016)..." + trailing tab
#      b) "Before analysis can begin... become skewed."
#      c) "After cleaning and schematizing... training times."
#      d) tab + "Another critical challenge..." + citation sdt for
#         Rawal et al. (2017) + closing sentence + bookmark.
# ---------------------------------------------------------------
$find.Execute("Before analysis can begin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$etlPara = $find.Parent.Paragraphs(1)
$etlRange = $etlPara.Range
$xmlEtl = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">According to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Gibert</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> et al. (2016), nearly 70% of all data mining occurs during the cleaning phase.</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Before analysis can begin, the data analyst needs to normalize the incoming data through an extract-transform-load (ETL) process.  This process needs to perform column renaming and reordering, adjusting quantity units, filtering erroneous values, populating missing values, and similar cleanup actions.  When analysis does not handle these aspects upfront, it creates a garbage-in/garbage-out scenario.  For example, a temperature reading of 55 degrees could be manually entered as 555 degrees, causing later analysis to become skewed.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>After cleaning and schematizing the incoming data, the next analyst needs to determine which aspects are relevant for their data mining objective.  Having large amounts of unrelated information does not improve results, and for many scenarios, it only slows down model training times.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>Another critical challenge is handling missing values</w:t></w:r><w:sdt><w:sdtPr><w:id w:val="1175543278"/><w:citation/></w:sdtPr><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Raw17 \l 1033 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> (Rawal et al., 2017)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> as they need to be normalized or removed.  These decisions become scenarios specific.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$etlRange.InsertXML($xmlEtl)

Write-Host "step4 done"

# ---------------------------------------------------------------
# 5) "Logical Component 2 of 2" gains a lastRenderedPageBreak marker
#    in front of its text run (another knock-on pagination shift).
# ---------------------------------------------------------------
$find.Execute("Logical Component 2 of 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$logicalPara = $find.Parent.Paragraphs(1)
$lr = $logicalPara.Range
$xmlLogical = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Logical Component 2 of 2</w:t></w:r></w:p>'
$lr.InsertXML($xmlLogical)

Write-Host "step5 done"
